# Sync non-localizable rule rows on the "Rules" sheet.
# The "BannedPaths" rule (row 35) is renamed to "BannedPath" and re-sorted to
# the bottom of this block (row 40, now with Severity "Critical" and no Tags),
# and the rows that used to follow it (36-40) each shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 <- old row 36 (CloudServiceIncompatibleWorkflowProcess)
$ws.Range("A35").Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Range("B35").Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Range("C35").Value = "Bug"
$ws.Range("D35").Value = "Blocker"
$ws.Range("E35").Value = "aem,cloud-service-compatibility"

# Row 36 <- old row 37 (IndexType)
$ws.Range("A36").Value = "IndexType"
$ws.Range("B36").Value = "Custom Search Index Definition Nodes Must Use the Index Type lucene"
$ws.Range("C36").Value = "Bug"
$ws.Range("D36").Value = "Blocker"
$ws.Range("E36").Value = "aem,cloud-service-compatibility"

# Row 37 <- old row 38 (IndexAsyncProperty)
$ws.Range("A37").Value = "IndexAsyncProperty"
$ws.Range("B37").Value = "Custom Lucene Oak Indexes must not be synchronous"
$ws.Range("C37").Value = "Bug"
$ws.Range("D37").Value = "Blocker"
$ws.Range("E37").Value = "aem,cloud-service-compatibility"

# Row 38 <- old row 39 (IndexTikaNode)
$ws.Range("A38").Value = "IndexTikaNode"
$ws.Range("B38").Value = "Custom Oak indexes must have a tika configuration"
$ws.Range("C38").Value = "Bug"
$ws.Range("D38").Value = "Blocker"
$ws.Range("E38").Value = "aem,cloud-service-compatibility"

# Row 39 <- old row 40 (IndexDamAssetLucene)
$ws.Range("A39").Value = "IndexDamAssetLucene"
$ws.Range("B39").Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Range("C39").Value = "Bug"
$ws.Range("D39").Value = "Blocker"
$ws.Range("E39").Value = "aem,cloud-service-compatibility"

# Row 40 <- renamed BannedPaths -> BannedPath, bumped to Critical, Tags cleared
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"
$ws.Range("E40").ClearContents()

# Update the saved selection to match the authored state.
$ws.Range("A37").Select() | Out-Null
